$d = $word.ActiveDocument

# --- New paragraph 1: "... i nomi e i cognomi dei clienti che hanno prenotazioni con arrivo precedente al primo luglio 2015"
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara1.Range.InsertAfter("Con riferimento al database IFTS scrivi la query che trova i nomi e i cognomi dei clienti che hanno prenotazioni con arrivo precedente al primo luglio 2015")

# --- New paragraph 2: "... città di residenza, nome e cognome dei clienti che hanno prenotazioni con caparra inferiore a 50 e importo superiore a 150"
$newPara1Again = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara1Again.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara2.Range.InsertAfter("Con riferimento al database IFTS scrivi la query che trova città di residenza, nome e cognome dei clienti che hanno prenotazioni con caparra inferiore a 50 e importo superiore a 150")

# --- New paragraph 3: empty paragraph, keeps the "Paragrafoelenco" style but no list numbering
$newPara2Again = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara2Again.Range.InsertParagraphAfter()
$newPara3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara3.Range.ParagraphFormat.Reset()
$newPara3.Style = "Paragrafoelenco"
